# Apply the edits described in the diff:
#  - Row 4 (person record) changes from Doug Derrick / Timepath Inc. / Analyst /
#    99 Shire Oak Road / dderrick@timepath.co.uk / 40799885412
#    to John Smith / IT Solutions / Analyst / 98 North Road /
#    jsmith@itsolutions.co.uk / 40716543298
#  - Row 10 (person record) changes from Jessie Marlowe / Aperture Inc. / Scientist /
#    27 Cheshire Street / jmarlowe@aperture.us / 40733154268
#    to Michelle Norton / Aperture Inc. / Scientist / 13 White Rabbit Street /
#    mnorton@aperture.us / 40731254562

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4
$ws.Range("A4").Value = "John"
$ws.Range("B4").Value = "Smith"
$ws.Range("C4").Value = "IT Solutions"
$ws.Range("D4").Value = "Analyst"
$ws.Range("E4").Value = "98 North Road"
$ws.Range("F4").Value = "jsmith@itsolutions.co.uk"
$ws.Range("G4").Value = 40716543298

# Update row 10
$ws.Range("A10").Value = "Michelle"
$ws.Range("B10").Value = "Norton"
$ws.Range("C10").Value = "Aperture Inc."
$ws.Range("D10").Value = "Scientist"
$ws.Range("E10").Value = "13 White Rabbit Street"
$ws.Range("F10").Value = "mnorton@aperture.us"
$ws.Range("G10").Value = 40731254562
